# Stundenliste-MichaelSchneider.xlsx - add new row for latest work session
# (Klassendiagramm besprochen, StateChart, Arbeitsaufteilung, loguru in App
# integriert. NEVER USE WINDOWS ...) plus adjust the running-total shared
# formula and the current view position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- copy formatting (number format / wrap text) from the previous row ----
$ws.Range("A27:D27").Copy()
$ws.Range("A28:D28").PasteSpecial(-4122)   # xlPasteFormats

# --- fill in the new row's data --------------------------------------------
$ws.Range("A28").Value2 = 44224
$ws.Range("B28").Value2 = 7
$ws.Range("C28").Formula = "=C27+B28"

$text = "Klassendiagramm besprochen, StateChart, Arbeitsaufteilung, loguru in App integriert.`nNEVER USE WINDOWS (daraus resultierte Fehler gefixt)"
$ws.Range("D28").Value = $text

# row 28 needs the same "tall" row height as the other multi-line entries
$ws.Rows.Item(28).RowHeight = 60

# --- move the current selection / view to the newly added row -------------
[void]$ws.Range("D29").Select()
